$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# --- 1. Mark two more chapters/sections as "Co" (Complete) -------------
$ws.Range("BN5").Value = "Co"
$ws.Range("BN6").Value = "Co"

# --- 2. Narrow the "status marker" columns (every other column) --------
# These columns hold the Status values (IP / Co) and were resized from an
# auto bestFit width down to a fixed narrow width.
$markerCols = @(
    "B","D","F","H","J","L","N","P","R","T","V","X","Z","AB","AD","AF",
    "AH","AJ","AL","AN","AP","AR","AT","AV","AX","AZ","BB","BD","BF","BH",
    "BJ","BL","BN","BP","BR","BT","BV","BX","BZ","CB","CD","CF","CH","CJ",
    "CL","CN"
)
foreach ($col in $markerCols) {
    $ws.Columns($col).ColumnWidth = 1.86
}

# --- 3. Extend the summary rows (23-25) so that every "data" column --
#        (not just the chapter/book totals AX:CN) is counted too.
$dataCols = @(
    "B","D","F","H","J","L","N","P","R","T","V","X","Z","AB","AD","AF",
    "AH","AJ","AL","AN","AP","AR","AT","AV"
)
foreach ($col in $dataCols) {
    $ws.Range($col + "23").Formula = "=COUNTIF(" + $col + "`$3:" + $col + "`$22,`$CO23)"
    $ws.Range($col + "24").Formula = "=COUNTIF(" + $col + "`$3:" + $col + "`$22,`$CO24)"
    $ws.Range($col + "25").Formula = "=" + $col + "`$2"
}

# --- 4. Row 26 no longer needs the per-book ratio in AX ---------------
$ws.Range("AX26").Clear()

# --- 5. Restore the view: no frozen/left-scrolled pane, selection on --
#        the last (percentage) row instead of the stray CP27 cell.
$ws.Range("B26:CN26").Select()
